$d = $word.ActiveDocument

$replacements = @(
    @("2024-01-19 Friday", "2024-01-20 Saturday"),
    @("549×7=", "111×2="),
    @("539×6=", "732×4="),
    @("834×8=", "220×9="),
    @("174×4=", "760×7="),
    @("490×3=", "461×7="),
    @("295×5=", "462×8="),
    @("534×3=", "864×4="),
    @("260×7=", "614×8="),
    @("383×8=", "576×7="),
    @("346×3=", "327×6="),
    @("934×3=", "499×2="),
    @("472×3=", "213×6="),
    @("536×9=", "228×7="),
    @("211×4=", "777×3="),
    @("195×8=", "544×6="),
    @("268×7=", "763×6="),
    @("259×8=", "964×9="),
    @("682×9=", "717×2="),
    @("855×8=", "184×5="),
    @("642×5=", "702×7="),
    @("309×2=", "587×6="),
    @("923×4=", "102×7="),
    @("191×3=", "780×6="),
    @("701×6=", "989×5="),
    @("135×2=", "691×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
